$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("B1").AddComment("test")
$tf2 = $c.Shape.TextFrame2
$tr = $tf2.TextRange
Write-Host "TextRange:" $tr.GetType()
Write-Host "Text:" $tr.Text
$tr.Font.Name = "Tahoma"
$tr.Font.Size = 9
Write-Host "B1 font after:" $ws.Range("B1").Font.Name
